$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")

# Official Development Assistance (ODA) row - Latest update date refreshed
$ws.Range("D11").Value = "Updated 21/01/2021"

# Social Institutions and Gender Index row - Latest update date refreshed
$ws.Range("D12").Value = "December, 2019"

# Wittgenstein Centre Human Capital Data Explorer row - year refreshed (2018 -> 2019).
# The target value looks like a plain number ("2019"), but the source cell is a text
# cell (same as the other "2019" text cells in this column), so force text storage
# via NumberFormat "@" and restore the cell's original look (style) with a
# formats-only paste from a cell that already holds that same style.
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2019"
$ws.Range("D19").PasteSpecial(-4122)

# World Development Indicators row - Latest update date refreshed
$ws.Range("D22").Value = "Updated 16/12/2020"

# World Urbanization Prospects row - revision year refreshed in both the
# description text and the year column (2018 -> 2019)
$ws.Range("B25").Value = "World Urbanization Prospects: the 2019 revision, UN"

$ws.Range("D14").Copy() | Out-Null
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2019"
$ws.Range("D25").PasteSpecial(-4122)

$excel.CutCopyMode = $false
